# Applies the Sept 27 2023 edit:
#  - On sheet "FirstLevelApprover", the hyperlinks in column B (B2:B6) are
#    re-pointed from mailto:Bingo@1234 to mailto:Bingo@123456 (and the
#    previous B2:B3 merged hyperlink range is split into individual
#    per-cell hyperlinks).
#  - The cell text for B2:B6 changes from "Bingo@1234" to "Bingo@123456".
#  - The active/selected sheet moves from "ExpenseRequest" to
#    "FirstLevelApprover", with the new selection sitting on B6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FirstLevelApprover")

$newAddress = "mailto:Bingo@123456"
$newDisplay = "Bingo@123456"

# Collect the addresses of the existing hyperlinks that live in column B
# (B2:B3 combined, B4, B5, B6) so we can remove them first. We gather the
# addresses into plain strings before deleting anything, since mutating
# the Hyperlinks collection while iterating it is not safe.
$addrsToRemove = @()
foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr.StartsWith('$B$')) {
        $addrsToRemove += $addr
    }
}

foreach ($addr in $addrsToRemove) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.Delete()
            break
        }
    }
}

# Re-create the hyperlinks and cell text for B2:B6 individually, each
# pointing at the new mailto address.
for ($r = 2; $r -le 6; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value = $newDisplay
    $ws.Hyperlinks.Add($cell, $newAddress)
    $cell.Style = "Hyperlink"
}

# Move the active sheet / tab selection from ExpenseRequest to
# FirstLevelApprover, and select B6 there.
$ws.Activate()
$ws.Range("B6").Select()
